$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 22, shifting rows 22-42 down to 23-43
$ws.Rows.Item(22).Insert()

# Fill in the new row 22 with data
$ws.Cells.Item(22, 1).Value = 1
$ws.Cells.Item(22, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(22, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(22, 4).Value = "2023-06-02"
$ws.Cells.Item(22, 5).Value = 15
$ws.Cells.Item(22, 6).Value = 100112003
$ws.Cells.Item(22, 7).Value = "Ajo"
$ws.Cells.Item(22, 8).Value = "Chino"
$ws.Cells.Item(22, 9).Value = "Primera"
$ws.Cells.Item(22, 10).Value = 400
$ws.Cells.Item(22, 11).Value = 16000
$ws.Cells.Item(22, 12).Value = 17000
$ws.Cells.Item(22, 13).Value = 16500
$ws.Cells.Item(22, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(22, 15).Value = "China"
$ws.Cells.Item(22, 16).Value = 1650
$ws.Cells.Item(22, 17).Value = 10
$ws.Cells.Item(22, 18).Value = "Hortaliza"
